# Mise à jour de l'application : ajout d'une nouvelle journée d'entrainement
# (nouvelle colonne AJ) dans le tableau de présences.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelle date d'entrainement en AJ1 (numéro de série Excel -> 20/08/2025)
$ws.Cells.Item(1, 36).Value = 45897

# Statuts de présence des joueurs pour cette nouvelle journée (colonne AJ, lignes 2 à 27)
$values = @("P","P","P","P","P","P","P","P","RH","P","P","A","P","P","P","B","P","P","P","P","P","P","P","P","P","RH")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 36).Value = $values[$i]
}

# Reproduit la mise en forme de la colonne précédente (AI) sur la nouvelle colonne (AJ)
$ws.Range("AI1:AI27").Copy()
$ws.Range("AJ1:AJ27").PasteSpecial(-4122)

# Ajuste le volet figé et la sélection active comme dans la nouvelle version du fichier
$win = $excel.ActiveWindow
$win.ScrollColumn = 33
$ws.Range("AL23").Select()
